$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "站点" (site) header/value from column A to column J.
# Column A now holds a numeric value (0) instead of the site name.
$ws.Range("J1").Value = $ws.Range("A1").Value()
$ws.Range("J2").Value = $ws.Range("A2").Value()

$ws.Range("A1").ClearContents()
$ws.Range("A2").Value = 0

# Update the active selection to match the new state.
$ws.Range("D6").Select()
